$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H33").Value = 230.41667
$ws_ALC.Range("I33").Value = 134.14635
$ws_ALC.Range("J33").Value = 794.2857
$ws_ALC.Range("K33").Value = 134.14635
$ws_ALC.Range("L33").Value = 794.2857
$ws_ALC.Range("M33").Value = 94.85364999999999
$ws_ALC.Range("N33").Value = -1252.2857
$ws_ALC.Range("H64").Value = 3355
$ws_ALC.Range("I64").Value = 2800
$ws_ALC.Range("J64").Value = 3416.6667
$ws_ALC.Range("K64").Value = 2800
$ws_ALC.Range("L64").Value = 3416.6667
$ws_ALC.Range("M64").Value = -2552
$ws_ALC.Range("N64").Value = -3912.6667
$ws_ALC.Range("H67").Value = 3355
$ws_ALC.Range("I67").Value = 2800
$ws_ALC.Range("J67").Value = 3416.6667
$ws_ALC.Range("K67").Value = 2800
$ws_ALC.Range("L67").Value = 3416.6667
$ws_ALC.Range("M67").Value = -1942
$ws_ALC.Range("N67").Value = -5132.6667
$ws_ALC.Range("H137").Value = 4624.6113
$ws_ALC.Range("I137").Value = 1080.8334
$ws_ALC.Range("J137").Value = 8168.3887
$ws_ALC.Range("K137").Value = 3242.5002
$ws_ALC.Range("L137").Value = 24505.1661
$ws_ALC.Range("M137").Value = -692.5002
$ws_ALC.Range("N137").Value = -29605.1661
$ws_ARM.Range("H32").Value = 13732.954
$ws_ARM.Range("I32").Value = 12380.9
$ws_ARM.Range("K32").Value = 12380.9
$ws_ARM.Range("M32").Value = -12093.9
$ws_ARM.Range("H45").Value = 1873.2858
$ws_ARM.Range("I45").Value = 1637.4546
$ws_ARM.Range("J45").Value = 2738
$ws_ARM.Range("K45").Value = 1637.4546
$ws_ARM.Range("L45").Value = 2738
$ws_ARM.Range("M45").Value = -1260.4546
$ws_ARM.Range("N45").Value = -3492
$ws_ARM.Range("H63").Value = 3411.25
$ws_ARM.Range("I63").Value = 2770
$ws_ARM.Range("K63").Value = 2770
$ws_ARM.Range("M63").Value = -2084
$ws_ARM.Range("H66").Value = 3411.25
$ws_ARM.Range("I66").Value = 2770
$ws_ARM.Range("K66").Value = 13850
$ws_ARM.Range("M66").Value = -10418
$ws_BSM.Range("H20").Value = 1548.1666
$ws_BSM.Range("I20").Value = 1078.6154
$ws_BSM.Range("K20").Value = 1078.6154
$ws_BSM.Range("M20").Value = -831.6153999999999
$ws_BSM.Range("H80").Value = 3939.4075
$ws_BSM.Range("I80").Value = 25103.25
$ws_BSM.Range("J80").Value = 258.73914
$ws_BSM.Range("K80").Value = 25103.25
$ws_BSM.Range("L80").Value = 258.73914
$ws_BSM.Range("M80").Value = -24105.25
$ws_BSM.Range("N80").Value = -2254.73914
$ws_BSM.Range("H83").Value = 3939.4075
$ws_BSM.Range("I83").Value = 25103.25
$ws_BSM.Range("J83").Value = 258.73914
$ws_BSM.Range("K83").Value = 125516.25
$ws_BSM.Range("L83").Value = 1293.6957
$ws_BSM.Range("M83").Value = -120524.25
$ws_BSM.Range("N83").Value = -11277.6957
$ws_BSM.Range("H94").Value = 1356.9032
$ws_BSM.Range("I94").Value = 1295.1428
$ws_BSM.Range("J94").Value = 1933.3334
$ws_BSM.Range("K94").Value = 1295.1428
$ws_BSM.Range("L94").Value = 1933.3334
$ws_BSM.Range("M94").Value = -844.1428000000001
$ws_BSM.Range("N94").Value = -2835.3334
$ws_BSM.Range("H105").Value = 2378.125
$ws_BSM.Range("I105").Value = 1701.6957
$ws_BSM.Range("K105").Value = 1701.6957
$ws_BSM.Range("M105").Value = 45.30430000000001
$ws_CRP.Range("H22").Value = 1705.6875
$ws_CRP.Range("I22").Value = 453.72726
$ws_CRP.Range("J22").Value = 4460
$ws_CRP.Range("K22").Value = 453.72726
$ws_CRP.Range("L22").Value = 4460
$ws_CRP.Range("M22").Value = -103.72726
$ws_CRP.Range("N22").Value = -5160
$ws_CRP.Range("H31").Value = 2836.25
$ws_CRP.Range("I31").Value = 795.76
$ws_CRP.Range("J31").Value = 3516.4133
$ws_CRP.Range("K31").Value = 795.76
$ws_CRP.Range("L31").Value = 3516.4133
$ws_CRP.Range("M31").Value = -500.76
$ws_CRP.Range("N31").Value = -4106.4133
$ws_CRP.Range("H34").Value = 2836.25
$ws_CRP.Range("I34").Value = 795.76
$ws_CRP.Range("J34").Value = 3516.4133
$ws_CRP.Range("K34").Value = 795.76
$ws_CRP.Range("L34").Value = 3516.4133
$ws_CRP.Range("M34").Value = -593.76
$ws_CRP.Range("N34").Value = -3920.4133
$ws_CRP.Range("H99").Value = 3084
$ws_CRP.Range("I99").Value = 2146.2856
$ws_CRP.Range("J99").Value = 5272
$ws_CRP.Range("K99").Value = 2146.2856
$ws_CRP.Range("L99").Value = 5272
$ws_CRP.Range("M99").Value = -648.2856000000002
$ws_CRP.Range("N99").Value = -8268
$ws_CRP.Range("H105").Value = 2483.45
$ws_CRP.Range("I105").Value = 2273.2778
$ws_CRP.Range("J105").Value = 4375
$ws_CRP.Range("K105").Value = 2273.2778
$ws_CRP.Range("L105").Value = 4375
$ws_CRP.Range("M105").Value = -526.2777999999998
$ws_CRP.Range("N105").Value = -7869
$ws_CRP.Range("H106").Value = 44488.5
$ws_CRP.Range("J106").Value = 44488.5
$ws_CRP.Range("L106").Value = 44488.5
$ws_CRP.Range("N106").Value = -47012.5
$ws_CRP.Range("H126").Value = 3084
$ws_CRP.Range("I126").Value = 2146.2856
$ws_CRP.Range("J126").Value = 5272
$ws_CRP.Range("K126").Value = 6438.8568
$ws_CRP.Range("L126").Value = 15816
$ws_CRP.Range("M126").Value = -3968.8568
$ws_CRP.Range("N126").Value = -20756
$ws_CUL.Range("H50").Value = 570.75
$ws_CUL.Range("I50").Value = 95
$ws_CUL.Range("J50").Value = 729.3333
$ws_CUL.Range("K50").Value = 285
$ws_CUL.Range("L50").Value = 2187.9999
$ws_CUL.Range("M50").Value = 196
$ws_CUL.Range("N50").Value = -3149.9999
$ws_CUL.Range("H53").Value = 570.75
$ws_CUL.Range("I53").Value = 95
$ws_CUL.Range("J53").Value = 729.3333
$ws_CUL.Range("K53").Value = 285
$ws_CUL.Range("L53").Value = 2187.9999
$ws_CUL.Range("M53").Value = 196
$ws_CUL.Range("N53").Value = -3149.9999
$ws_CUL.Range("H69").Value = 44758880
$ws_CUL.Range("I69").Value = 421.66666
$ws_CUL.Range("J69").Value = 60555980
$ws_CUL.Range("K69").Value = 1264.99998
$ws_CUL.Range("L69").Value = 181667940
$ws_CUL.Range("M69").Value = -453.9999800000001
$ws_CUL.Range("N69").Value = -181669562
$ws_CUL.Range("H72").Value = 44758880
$ws_CUL.Range("I72").Value = 421.66666
$ws_CUL.Range("J72").Value = 60555980
$ws_CUL.Range("K72").Value = 3794.99994
$ws_CUL.Range("L72").Value = 545003820
$ws_CUL.Range("M72").Value = 261.0000600000003
$ws_CUL.Range("N72").Value = -545011932
$ws_CUL.Range("H80").Value = 71572424
$ws_CUL.Range("I80").Value = 250874.88
$ws_CUL.Range("J80").Value = 166667820
$ws_CUL.Range("K80").Value = 752624.64
$ws_CUL.Range("L80").Value = 500003460
$ws_CUL.Range("M80").Value = -751688.64
$ws_CUL.Range("N80").Value = -500005332
$ws_CUL.Range("H83").Value = 71572424
$ws_CUL.Range("I83").Value = 250874.88
$ws_CUL.Range("J83").Value = 166667820
$ws_CUL.Range("K83").Value = 2257873.92
$ws_CUL.Range("L83").Value = 1500010380
$ws_CUL.Range("M83").Value = -2253193.92
$ws_CUL.Range("N83").Value = -1500019740
$ws_CUL.Range("H102").Value = 19800
$ws_CUL.Range("J102").Value = 19800
$ws_CUL.Range("L102").Value = 59400
$ws_CUL.Range("N102").Value = -64268
$ws_CUL.Range("H113").Value = 3779.2424
$ws_CUL.Range("I113").Value = 6435.1177
$ws_CUL.Range("J113").Value = 957.375
$ws_CUL.Range("K113").Value = 19305.3531
$ws_CUL.Range("L113").Value = 2872.125
$ws_CUL.Range("M113").Value = -17135.3531
$ws_CUL.Range("N113").Value = -7212.125
$ws_CUL.Range("H129").Value = 126435.375
$ws_CUL.Range("I129").Value = 300737.7
$ws_CUL.Range("J129").Value = 1933.7142
$ws_CUL.Range("K129").Value = 902213.1000000001
$ws_CUL.Range("L129").Value = 5801.142599999999
$ws_CUL.Range("M129").Value = -897213.1000000001
$ws_CUL.Range("N129").Value = -15801.1426
$ws_CUL.Range("H132").Value = 2328.7407
$ws_CUL.Range("I132").Value = 1676.1538
$ws_CUL.Range("J132").Value = 2934.7144
$ws_CUL.Range("K132").Value = 15085.3842
$ws_CUL.Range("L132").Value = 26412.4296
$ws_CUL.Range("M132").Value = -12555.3842
$ws_CUL.Range("N132").Value = -31472.4296
$ws_GSM.Range("H80").Value = 5658.6
$ws_GSM.Range("J80").Value = 4655.143
$ws_GSM.Range("L80").Value = 4655.143
$ws_GSM.Range("N80").Value = -6651.143
$ws_GSM.Range("H83").Value = 5658.6
$ws_GSM.Range("J83").Value = 4655.143
$ws_GSM.Range("L83").Value = 23275.715
$ws_GSM.Range("N83").Value = -33259.715
$ws_GSM.Range("H113").Value = 1660.6923
$ws_GSM.Range("I113").Value = 1568.9
$ws_GSM.Range("J113").Value = 1966.6666
$ws_GSM.Range("K113").Value = 1568.9
$ws_GSM.Range("L113").Value = 1966.6666
$ws_GSM.Range("M113").Value = 601.0999999999999
$ws_GSM.Range("N113").Value = -6306.6666
$ws_GSM.Range("H132").Value = 3334.9565
$ws_GSM.Range("I132").Value = 2399.75
$ws_GSM.Range("J132").Value = 3833.7334
$ws_GSM.Range("K132").Value = 7199.25
$ws_GSM.Range("L132").Value = 11501.2002
$ws_GSM.Range("M132").Value = -4669.25
$ws_GSM.Range("N132").Value = -16561.2002
$ws_LTW.Range("H16").Value = 1576.4412
$ws_LTW.Range("I16").Value = 1460
$ws_LTW.Range("J16").Value = 2449.75
$ws_LTW.Range("K16").Value = 1460
$ws_LTW.Range("L16").Value = 2449.75
$ws_LTW.Range("M16").Value = -1290
$ws_LTW.Range("N16").Value = -2789.75
$ws_LTW.Range("H93").Value = 1520.3077
$ws_LTW.Range("I93").Value = 1049.0769
$ws_LTW.Range("J93").Value = 1991.5385
$ws_LTW.Range("K93").Value = 1049.0769
$ws_LTW.Range("L93").Value = 1991.5385
$ws_LTW.Range("M93").Value = 198.9231
$ws_LTW.Range("N93").Value = -4487.538500000001
$ws_LTW.Range("H132").Value = 2808.1064
$ws_LTW.Range("I132").Value = 2149.6667
$ws_LTW.Range("K132").Value = 6449.000100000001
$ws_LTW.Range("M132").Value = -3919.000100000001
$ws_WVR.Range("H81").Value = 3132.6
$ws_WVR.Range("I81").Value = 1013.63635
$ws_WVR.Range("J81").Value = 5722.4443
$ws_WVR.Range("K81").Value = 2027.2727
$ws_WVR.Range("L81").Value = 11444.8886
$ws_WVR.Range("M81").Value = -966.2727
$ws_WVR.Range("N81").Value = -13566.8886
$ws_WVR.Range("H84").Value = 3132.6
$ws_WVR.Range("I84").Value = 1013.63635
$ws_WVR.Range("J84").Value = 5722.4443
$ws_WVR.Range("K84").Value = 10136.3635
$ws_WVR.Range("L84").Value = 57224.443
$ws_WVR.Range("M84").Value = -4832.363499999999
$ws_WVR.Range("N84").Value = -67832.443
$ws_WVR.Range("H96").Value = 1300
$ws_WVR.Range("I96").Value = 0
$ws_WVR.Range("J96").Value = 1300
$ws_WVR.Range("K96").Value = 0
$ws_WVR.Range("L96").Value = 1300
$ws_WVR.Range("M96").ClearContents()
$ws_WVR.Range("N96").Value = -4046
